$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 15, shifting existing rows 15-26 down to 16-27.
$ws.Rows.Item(15).Insert()

# Populate the newly inserted row 15 with the new data record.
$ws.Cells.Item(15, 1).Value = 6
$ws.Cells.Item(15, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(15, 3).Value = "Metropolitana"
$ws.Cells.Item(15, 4).Value = 44755
$ws.Cells.Item(15, 5).Value = 13
$ws.Cells.Item(15, 6).Value = 100112035
$ws.Cells.Item(15, 7).Value = "Bruselas (repollito)"
$ws.Cells.Item(15, 8).Value = "Sin especificar"
$ws.Cells.Item(15, 9).Value = "Primera"
$ws.Cells.Item(15, 10).Value = 230
$ws.Cells.Item(15, 11).Value = 16000
$ws.Cells.Item(15, 12).Value = 18000
$ws.Cells.Item(15, 13).Value = 16783
$ws.Cells.Item(15, 14).Value = "`$/malla 15 kilos"
$ws.Cells.Item(15, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(15, 16).Value = 1119
$ws.Cells.Item(15, 17).Value = 15
$ws.Cells.Item(15, 18).Value = "Hortaliza"

Write-Output "done"
